$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 382/383, shifting the existing rows 382-400 down to 384-402
$ws.Rows("382:383").Insert()

# New row 382 - latest weekly price entry
$ws.Range("A382").Value = 1
$ws.Range("B382").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C382").Value = "Arica y Parinacota"
$ws.Range("D382").Value = 44939
$ws.Range("E382").Value = 15
$ws.Range("F382").Value = 100114013
$ws.Range("G382").Value = "Zanahoria"
$ws.Range("H382").Value = "Sin especificar"
$ws.Range("I382").Value = "Primera"
$ws.Range("J382").Value = 100
$ws.Range("K382").Value = 16000
$ws.Range("L382").Value = 17000
$ws.Range("M382").Value = 16600
$ws.Range("N382").Value = "$/saco 25 kilos"
$ws.Range("O382").Value = "Región de Arica y Parinacota"
$ws.Range("P382").Value = 664
$ws.Range("Q382").Value = 25
$ws.Range("R382").Value = "Hortaliza"

# New row 383 - latest weekly price entry
$ws.Range("A383").Value = 1
$ws.Range("B383").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C383").Value = "Arica y Parinacota"
$ws.Range("D383").Value = 44939
$ws.Range("E383").Value = 15
$ws.Range("F383").Value = 100114013
$ws.Range("G383").Value = "Zanahoria"
$ws.Range("H383").Value = "Sin especificar"
$ws.Range("I383").Value = "Primera"
$ws.Range("J383").Value = 45
$ws.Range("K383").Value = 16000
$ws.Range("L383").Value = 17000
$ws.Range("M383").Value = 16556
$ws.Range("N383").Value = "$/saco 25 kilos"
$ws.Range("O383").Value = "Valle de Camiña"
$ws.Range("P383").Value = 662
$ws.Range("Q383").Value = 25
$ws.Range("R383").Value = "Hortaliza"
